$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values of 45185 (2023-09-16) are updated to 45204 (2023-10-05)
# for all data rows (rows 2 through 70).
$ws.Range("C2:C70").Value = 45204
